# Teachers free time variable added:
#  - Remove the unused "Saturday" column from the TeacherFreeSlot sheet
#  - Re-order the teacher rows (AI above MMR)
#  - Update the view selection accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TeacherFreeSlot")
$ws.Activate()

# Remove the (empty) Saturday column; this shifts Sunday..Thursday one column to the left
# and keeps each teacher's day values aligned with their original row.
$ws.Columns("B").Delete()

# Swap the teacher-initial labels so that AI's row comes before MMR's row,
# while each row keeps its own free-time values.
$teacherRow2 = $ws.Cells.Item(2, 1).Value
$teacherRow3 = $ws.Cells.Item(3, 1).Value
$ws.Cells.Item(2, 1).Value = $teacherRow3
$ws.Cells.Item(3, 1).Value = $teacherRow2

# Update the selection / view to match the edited sheet
$ws.Range("F1").Select()
